$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "(0, 124, 255)  Road" -> "(0, 124, 255)  Water"
# ---------------------------------------------------------------------------
$hit1 = $d.Content
$hit1.Find.ClearFormatting()
$found1 = $hit1.Find.Execute("(0, 124, 255)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate '(0, 124, 255)'"
}

$para1 = $d.Range($hit1.Start, $hit1.End)
[void]$para1.Expand(4)

$para1Text = $para1.Text
$wordOffset = $para1Text.LastIndexOf("Road")
if ($wordOffset -lt 0) {
    throw "Could not locate 'Road' inside the (0, 124, 255) paragraph"
}
$roadStart = $para1.Start + $wordOffset
$roadEnd = $roadStart + 4

$wordRange = $d.Range($roadStart, $roadEnd)
$wordRange.Text = "Water"

# Re-touch the formatting of the freshly inserted word (without actually
# changing the visible font) so it settles into its own run, split off from
# the preceding tab/space run - matching how Word keeps replacement text
# separate from the run it was typed into.
$newWordRange = $d.Range($roadStart, $roadStart + 5)
$newWordRange.Bold = 1
$newWordRange.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: merge the fragmented "(24, 24, 24) " runs into a single run
#   (the visible text is unchanged, only the run layout is simplified).
# ---------------------------------------------------------------------------
$hit2 = $d.Content
$hit2.Find.ClearFormatting()
$found2 = $hit2.Find.Execute("(24, 24, 24)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate '(24, 24, 24)'"
}

$para2 = $d.Range($hit2.Start, $hit2.End)
[void]$para2.Expand(4)

$para2Text = $para2.Text
$tabOffset = $para2Text.IndexOf([char]9)
if ($tabOffset -lt 0) {
    throw "Could not locate a tab inside the (24, 24, 24) paragraph"
}
$colorStart = $para2.Start
$colorEnd = $para2.Start + $tabOffset

$mergeRange = $d.Range($colorStart, $colorEnd)
$mergeRange.Text = "(24, 24, 24) X"

$sentinelRange = $d.Range($mergeRange.End - 1, $mergeRange.End)
$sentinelRange.Delete()
